$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Service Contacts sheet: move "delivery_organisation_path" (was column R)
# to become the new column D, immediately before "practitioner_key"
# (which, together with everything from the old practitioner_key column
# through funding_source, shifts one column to the right: D..Q -> E..R).
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Service Contacts")

# Row 1 (headers)
$ws.Range("D1").Value = "delivery_organisation_path"
$ws.Range("E1").Value = "practitioner_key"
$ws.Range("F1").Value = "service_contact_date"
$ws.Range("G1").Value = "service_contact_type"
$ws.Range("H1").Value = "service_contact_postcode"
$ws.Range("I1").Value = "service_contact_modality"
$ws.Range("J1").Value = "service_contact_participants"
$ws.Range("K1").Value = "service_contact_venue"
$ws.Range("L1").Value = "service_contact_duration"
$ws.Range("M1").Value = "service_contact_copayment"
$ws.Range("N1").Value = "service_contact_participation_indicator"
$ws.Range("O1").Value = "service_contact_interpreter"
$ws.Range("P1").Value = "service_contact_no_show"
$ws.Range("Q1").Value = "service_contact_final"
$ws.Range("R1").Value = "funding_source"

# Row 2 (CL0001-E01-SC01)
$ws.Range("D2").Value = "PHN999:NFP02"
$ws.Range("E2").Value = "P01"
$ws.Range("F2").Value = 21052016
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 9999
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 1

# Row 3 (CL0002-E01-SC01)
$ws.Range("D3").Value = "PHN999:NFP01"
$ws.Range("D3").Font.Color = 0
$ws.Range("E3").Value = "P02"
$ws.Range("F3").Value = 15062016
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 9999
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 2

# Column widths: the wide (19-char) formatting that used to sit on
# P:R (service_contact_final, funding_source, delivery_organisation_path)
# now sits on D (delivery_organisation_path) and Q:R (service_contact_final,
# funding_source).
$ws.Columns.Item(4).ColumnWidth = 18.14
$ws.Range($ws.Columns.Item(17), $ws.Columns.Item(18)).ColumnWidth = 18.14

# Selection / view left on the new delivery_organisation_path column.
$ws.Range("D1:D1048576").Select() | Out-Null

# -----------------------------------------------------------------------
# Other sheets: selection-only changes (no data/format changes), left over
# from the author navigating the workbook while making the edit above.
# -----------------------------------------------------------------------
$wsOrg = $wb.Worksheets.Item("Organisations")
$wsOrg.Range("H1:J3").Select() | Out-Null

$wsK10 = $wb.Worksheets.Item("K10+")
$wsK10.Range("F1:F5").Select() | Out-Null

$wsK5 = $wb.Worksheets.Item("K5")
$wsK5.Range("F1:F5").Select() | Out-Null

# Leave the Service Contacts sheet active with its new selection, matching
# the saved workbook state.
$ws.Activate() | Out-Null
$ws.Range("D1:D1048576").Select() | Out-Null
